$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.180.02'
$ws.Range('E2').Value = '  -1.26%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.861.43'
$ws.Range('E3').Value = '  -1.03%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7125'
$ws.Range('E5').Value = '  -2.02%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '240.53'
$ws.Range('E6').Value = '  +0.47%  '

$ws.Range('E7').Value = '  +0.24%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07705'
$ws.Range('E8').Value = '  -2.21%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3072'
$ws.Range('E9').Value = '  -0.56%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.89'
$ws.Range('E10').Value = '  -1.71%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08261'
$ws.Range('E11').Value = '  +0.57%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.857.54'
$ws.Range('E12').Value = '  -0.94%  '

$ws.Range('E13').Value = '  -1.45%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.210'
$ws.Range('E14').Value = '  -1.20%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.20'
$ws.Range('E15').Value = '  +0.55%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.177.04'
$ws.Range('E16').Value = '  -1.24%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.856'
$ws.Range('E17').Value = '  -0.03%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '243.16'
$ws.Range('E18').Value = '  +0.37%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007795'
$ws.Range('E19').Value = '  -0.99%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.13'
$ws.Range('E20').Value = '  -1.92%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.107.88'
$ws.Range('E21').Value = '  -0.77%  '

$ws.Range('E22').Value = '  +0.19%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.937'
$ws.Range('E23').Value = '  +1.98%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.13%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1572'
$ws.Range('E25').Value = '  +5.23%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.40'
$ws.Range('E26').Value = '  -0.31%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.896'
$ws.Range('E27').Value = '  -1.20%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.18'
$ws.Range('E28').Value = '  -0.60%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.327'
$ws.Range('E29').Value = '  -2.54%  '

$ws.Range('E30').Value = '  +1.09%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.344'
$ws.Range('E31').Value = '  -0.48%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.085'
$ws.Range('E32').Value = '  -0.51%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05177'
$ws.Range('E33').Value = '  -1.46%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.902'
$ws.Range('E34').Value = '  -2.89%  '

$ws.Range('E35').Value = '  -2.19%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7265'
$ws.Range('E36').Value = '  +1.01%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.683'
$ws.Range('E37').Value = '  +0.46%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01846'
$ws.Range('E38').Value = '  -0.85%  '

$ws.Range('E39').Value = '  -1.14%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.140.70'
$ws.Range('E40').Value = '  -3.01%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8986'
$ws.Range('E41').Value = '  -1.45%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.076'
$ws.Range('E42').Value = '  +1.33%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.04'
$ws.Range('E43').Value = '  -0.01%  '

$ws.Range('E44').Value = '  +0.24%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.41'
$ws.Range('E45').Value = '  -1.15%  '

$ws.Range('E46').Value = '  -1.53%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.003.74'
$ws.Range('E47').Value = '  -0.80%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.759'
$ws.Range('E48').Value = '  -1.13%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.278'
$ws.Range('E49').Value = '  +0.13%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.857'
$ws.Range('E50').Value = '  -1.12%  '

$ws.Range('E51').Value = '  -0.31%  '
